# Auto-generated edit script applying the TMS_PCMs.xlsx ticket-data refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").NumberFormat = "@"
$ws.Range("B2").Value = "TT-20251024000052"
$ws.Range("F2").NumberFormat = "@"
$ws.Range("F2").Value = "2025-10-24 04:26:39"
$ws.Range("I2").NumberFormat = "@"
$ws.Range("I2").Value = "MKTF0647"
$ws.Range("J2").NumberFormat = "@"
$ws.Range("J2").Value = "Region_5"
$ws.Range("K2").NumberFormat = "@"
$ws.Range("K2").Value = "Generator_SG"
$ws.Range("L2").NumberFormat = "@"
$ws.Range("L2").Value = "PCM-20251024-00000228"
$ws.Range("O2").NumberFormat = "@"
$ws.Range("O2").Value = "1(2G/3G) sites down under MKTF0647-P1-USF/MKTF0647-P1-USF @ Generator_SG"
$ws.Range("P2").NumberFormat = "@"
$ws.Range("P2").Value = "P1"
$ws.Range("S2").NumberFormat = "@"
$ws.Range("S2").Value = "PCM-20251024-00000228"
$ws.Range("T2").NumberFormat = "@"
$ws.Range("T2").Value = "TT-20251024-00268"
$ws.Range("W2").NumberFormat = "@"
$ws.Range("W2").Value = "Critical"
$ws.Range("X2").NumberFormat = "@"
$ws.Range("X2").Value = "2025-10-24 04:18:33"
$ws.Range("AC2").NumberFormat = "@"
$ws.Range("AC2").Value = "2025-10-24 04:18:33"
$ws.Range("AI2").NumberFormat = "@"
$ws.Range("AI2").Value = "PCM-20251024-00000228"
$ws.Range("AK2").NumberFormat = "@"
$ws.Range("AK2").Value = "MKTF0647"
$ws.Range("AM2").NumberFormat = "@"
$ws.Range("AM2").Value = "makkah"
$ws.Range("AS2").NumberFormat = "@"
$ws.Range("AS2").Value = "MKTF0647-P1-USF"
$ws.Range("AT2").NumberFormat = "@"
$ws.Range("AT2").Value = "SERIAL-20251024-0052"
$ws.Range("AU2").NumberFormat = "@"
$ws.Range("AU2").Value = "2025-10-24 04:18:36"
$ws.Range("B3").NumberFormat = "@"
$ws.Range("B3").Value = "TT-20251024000050"
$ws.Range("F3").NumberFormat = "@"
$ws.Range("F3").Value = "2025-10-24 04:17:24"
$ws.Range("I3").NumberFormat = "@"
$ws.Range("I3").Value = "HMBY2212"
$ws.Range("J3").NumberFormat = "@"
$ws.Range("J3").Value = "Region_1"
$ws.Range("K3").NumberFormat = "@"
$ws.Range("K3").Value = "Generator_SG"
$ws.Range("L3").NumberFormat = "@"
$ws.Range("L3").Value = "PCM-20251024-00000218"
$ws.Range("O3").NumberFormat = "@"
$ws.Range("O3").Value = "2(2G) 1(3G) 2(LTE) sites down under HMBY2212-P3,KBR2492-P3/DM1MBSCH01 @ Sceco"
$ws.Range("S3").NumberFormat = "@"
$ws.Range("S3").Value = "PCM-20251024-00000218"
$ws.Range("T3").NumberFormat = "@"
$ws.Range("T3").Value = "TT-20251024-00248"
$ws.Range("V3").NumberFormat = "@"
$ws.Range("V3").Value = "CSL Fault"
$ws.Range("X3").NumberFormat = "@"
$ws.Range("X3").Value = "2025-10-24 04:07:15"
$ws.Range("Y3").NumberFormat = "@"
$ws.Range("Y3").Value = "2G_3G_LTE SA"
$ws.Range("Z3").NumberFormat = "@"
$ws.Range("Z3").Value = "21825"
$ws.Range("AA3").NumberFormat = "@"
$ws.Range("AA3").Value = "CSL Fault"
$ws.Range("AC3").NumberFormat = "@"
$ws.Range("AC3").Value = "2025-10-24 04:07:15"
$ws.Range("AE3").NumberFormat = "@"
$ws.Range("AE3").Value = "Dear Team,`nHub:- HMBY2212`nSuspecting power issue at hub site and kindly send EM-FLM team to check and resolve issue at site.`nKindly update TT under EM FLM"
$ws.Range("AF3").NumberFormat = "@"
$ws.Range("AF3").Value = "Dear Team,`nHub:- HMBY2212`nSuspecting power issue at hub site and kindly send EM-FLM team to check and resolve issue at site.`nKindly update TT under EM FLM"
$ws.Range("AG3").NumberFormat = "@"
$ws.Range("AG3").Value = "PCM-20251024-00000215"
$ws.Range("AI3").NumberFormat = "@"
$ws.Range("AI3").Value = "PCM-20251024-00000218"
$ws.Range("AK3").NumberFormat = "@"
$ws.Range("AK3").Value = "HMBY2212"
$ws.Range("AM3").NumberFormat = "@"
$ws.Range("AM3").Value = "dammam"
$ws.Range("AO3").NumberFormat = "@"
$ws.Range("AO3").Value = "Telecom"
$ws.Range("AP3").NumberFormat = "@"
$ws.Range("AP3").Value = "SRAN"
$ws.Range("AS3").NumberFormat = "@"
$ws.Range("AS3").Value = "KBR2492-P3 2G_LTE,HMBY2212-P3 2G_3G_LTE"
$ws.Range("AT3").NumberFormat = "@"
$ws.Range("AT3").Value = "SERIAL-20251024-0050"
$ws.Range("AU3").NumberFormat = "@"
$ws.Range("AU3").Value = "2025-10-24 04:07:17"
$ws.Range("B4").NumberFormat = "@"
$ws.Range("B4").Value = "TT-20251024000043"
$ws.Range("F4").NumberFormat = "@"
$ws.Range("F4").Value = "2025-10-24 03:23:39"
$ws.Range("I4").NumberFormat = "@"
$ws.Range("I4").Value = "HAW0293"
$ws.Range("J4").NumberFormat = "@"
$ws.Range("J4").Value = "Region_4"
$ws.Range("L4").NumberFormat = "@"
$ws.Range("L4").Value = "PCM-20251024-00000171"
$ws.Range("O4").NumberFormat = "@"
$ws.Range("O4").Value = "1(2G/3G/5G) sites down under HAW0293-P2/HAW0293-P2 @ Sceco"
$ws.Range("P4").NumberFormat = "@"
$ws.Range("P4").Value = "P2"
$ws.Range("S4").NumberFormat = "@"
$ws.Range("S4").Value = "PCM-20251024-00000171"
$ws.Range("T4").NumberFormat = "@"
$ws.Range("T4").Value = "TT-20251024-00200"
$ws.Range("W4").NumberFormat = "@"
$ws.Range("W4").Value = "Major"
$ws.Range("X4").NumberFormat = "@"
$ws.Range("X4").Value = "2025-10-24 03:15:23"
$ws.Range("Y4").NumberFormat = "@"
$ws.Range("Y4").Value = "2G_3G SA"
$ws.Range("AC4").NumberFormat = "@"
$ws.Range("AC4").Value = "2025-10-24 03:15:23"
$ws.Range("AE4").NumberFormat = "@"
$ws.Range("AE4").Value = ""
$ws.Range("AF4").NumberFormat = "@"
$ws.Range("AF4").Value = ""
$ws.Range("AG4").NumberFormat = "@"
$ws.Range("AG4").Value = ""
$ws.Range("AI4").NumberFormat = "@"
$ws.Range("AI4").Value = "PCM-20251024-00000171"
$ws.Range("AK4").NumberFormat = "@"
$ws.Range("AK4").Value = "HAW0293"
$ws.Range("AM4").NumberFormat = "@"
$ws.Range("AM4").Value = "taif"
$ws.Range("AS4").NumberFormat = "@"
$ws.Range("AS4").Value = "HAW0293-P2"
$ws.Range("AT4").NumberFormat = "@"
$ws.Range("AT4").Value = "SERIAL-20251024-0042"
$ws.Range("AU4").NumberFormat = "@"
$ws.Range("AU4").Value = "2025-10-24 03:15:26"
$ws.Range("B5").NumberFormat = "@"
$ws.Range("B5").Value = "TT-20251024000028"
$ws.Range("F5").NumberFormat = "@"
$ws.Range("F5").Value = "2025-10-24 02:23:59"
$ws.Range("I5").NumberFormat = "@"
$ws.Range("I5").Value = "MAK0195"
$ws.Range("J5").NumberFormat = "@"
$ws.Range("J5").Value = "Region_5"
$ws.Range("K5").NumberFormat = "@"
$ws.Range("K5").Value = "Sceco"
$ws.Range("L5").NumberFormat = "@"
$ws.Range("L5").Value = "PCM-20251024-00000121"
$ws.Range("O5").NumberFormat = "@"
$ws.Range("O5").Value = "1(2G/5G) sites down under MAK0195-P1-HUB/MAK0195-P1-HUB @ Sceco"
$ws.Range("P5").NumberFormat = "@"
$ws.Range("P5").Value = "P1"
$ws.Range("S5").NumberFormat = "@"
$ws.Range("S5").Value = "PCM-20251024-00000121"
$ws.Range("T5").NumberFormat = "@"
$ws.Range("T5").Value = "TT-20251024-00113"
$ws.Range("V5").NumberFormat = "@"
$ws.Range("V5").Value = "BTS O&M LINK FAILURE"
$ws.Range("W5").NumberFormat = "@"
$ws.Range("W5").Value = "Critical"
$ws.Range("X5").NumberFormat = "@"
$ws.Range("X5").Value = "2025-10-24 01:51:15"
$ws.Range("Y5").NumberFormat = "@"
$ws.Range("Y5").Value = "2G SA"
$ws.Range("Z5").NumberFormat = "@"
$ws.Range("Z5").Value = "7706"
$ws.Range("AA5").NumberFormat = "@"
$ws.Range("AA5").Value = "BTS O&M LINK FAILURE"
$ws.Range("AC5").NumberFormat = "@"
$ws.Range("AC5").Value = "2025-10-24 01:51:15"
$ws.Range("AE5").NumberFormat = "@"
$ws.Range("AE5").Value = "Power issue at site"
$ws.Range("AF5").NumberFormat = "@"
$ws.Range("AF5").Value = "Power issue at site"
$ws.Range("AG5").NumberFormat = "@"
$ws.Range("AG5").Value = "PCM-20251024-00000102"
$ws.Range("AI5").NumberFormat = "@"
$ws.Range("AI5").Value = "PCM-20251024-00000121"
$ws.Range("AK5").NumberFormat = "@"
$ws.Range("AK5").Value = "MAK0195"
$ws.Range("AM5").NumberFormat = "@"
$ws.Range("AM5").Value = "makkah"
$ws.Range("AS5").NumberFormat = "@"
$ws.Range("AS5").Value = "MAK0195-P1-HUB"
$ws.Range("AT5").NumberFormat = "@"
$ws.Range("AT5").Value = "SERIAL-20251024-0028"
$ws.Range("AU5").NumberFormat = "@"
$ws.Range("AU5").Value = "2025-10-24 01:51:23"
$ws.Range("B6").NumberFormat = "@"
$ws.Range("B6").Value = "TT-20251024000023"
$ws.Range("F6").NumberFormat = "@"
$ws.Range("F6").Value = "2025-10-24 02:09:09"
$ws.Range("I6").NumberFormat = "@"
$ws.Range("I6").Value = "HMBY0419"
$ws.Range("J6").NumberFormat = "@"
$ws.Range("J6").Value = "Region_1"
$ws.Range("K6").NumberFormat = "@"
$ws.Range("K6").Value = "Generator_SG"
$ws.Range("L6").NumberFormat = "@"
$ws.Range("L6").Value = "PCM-20251024-00000107"
$ws.Range("O6").NumberFormat = "@"
$ws.Range("O6").Value = "1(2G/LTE) sites down under HMBY0419-P3/DM1MBSCH01 @ Generator_SG"
$ws.Range("P6").NumberFormat = "@"
$ws.Range("P6").Value = "P3"
$ws.Range("S6").NumberFormat = "@"
$ws.Range("S6").Value = "PCM-20251024-00000107"
$ws.Range("T6").NumberFormat = "@"
$ws.Range("T6").Value = "TT-20251024-00121"
$ws.Range("W6").NumberFormat = "@"
$ws.Range("W6").Value = "Major"
$ws.Range("X6").NumberFormat = "@"
$ws.Range("X6").Value = "2025-10-24 02:00:47"
$ws.Range("AC6").NumberFormat = "@"
$ws.Range("AC6").Value = "2025-10-24 02:00:47"
$ws.Range("AE6").NumberFormat = "@"
$ws.Range("AE6").Value = ""
$ws.Range("AF6").NumberFormat = "@"
$ws.Range("AF6").Value = ""
$ws.Range("AG6").NumberFormat = "@"
$ws.Range("AG6").Value = ""
$ws.Range("AI6").NumberFormat = "@"
$ws.Range("AI6").Value = "PCM-20251024-00000107"
$ws.Range("AK6").NumberFormat = "@"
$ws.Range("AK6").Value = "HMBY0419"
$ws.Range("AM6").NumberFormat = "@"
$ws.Range("AM6").Value = "dammam"
$ws.Range("AO6").NumberFormat = "@"
$ws.Range("AO6").Value = "EM"
$ws.Range("AP6").NumberFormat = "@"
$ws.Range("AP6").Value = "2G"
$ws.Range("AS6").NumberFormat = "@"
$ws.Range("AS6").Value = "HMBY0419-P3"
$ws.Range("AT6").NumberFormat = "@"
$ws.Range("AT6").Value = "SERIAL-20251024-0023"
$ws.Range("AU6").NumberFormat = "@"
$ws.Range("AU6").Value = "2025-10-24 02:00:49"
$ws.Range("M8").NumberFormat = "@"
$ws.Range("M8").Value = "Warning"

Write-Host "Applied all cell updates"
